$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")
$ws.Range("F16").Value = "10/10/2019"
